# Apply updated TPM-derived NATMI ligand-receptor statistics to Sheet1 (rows 2-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('ECs', 'Fgf17', 'Fgfr4', 'ECs', 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 3, 1, 0.903185, 2.709555, 0.03154869388788047, 0.03154869388788046, 0.1640900486766667, 1.47681043809, 0.0111832175479874, 0.0111832175479874),
    @('ECs', 'Fgf17', 'Fgfr4', 'FAPs', 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 2, 0.6666666666666666, 0.110372, 0.331116, 0.003855347953955327, 0.003855347953955326, 0.02005231137866666, 0.180470802408, 0.001366623767230927, 0.001366623767230927),
    @('ECs', 'Fgf17', 'Fgfr4', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 1, 0.3333333333333333, 0.03076233333333333, 0.09228699999999999, 0.001074543352259254, 0.001074543352259254, 0.005588880211777778, 0.050299921906, 0.0003808985600407125, 0.0003808985600407127),
    @('ECs', 'Fgf17', 'Fgfr4', 'MuSCs', 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 3, 1, 27.52907633333334, 82.58722900000001, 0.9616041035407232, 0.9616041035407231, 5.001464235522445, 45.013178119702, 0.3408644403204415, 0.3408644403204416),
    @('ECs', 'Fgf17', 'Fgfr4', 'Resolving-Mac', 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 2, 0.6666666666666666, 0.05488933333333334, 0.164668, 0.001917311265181737, 0.001917311265181736, 0.009972257487111112, 0.08975031738400001, 0.0006796385632297515, 0.0006796385632297516),
    @('FAPs', 'Fgf17', 'Fgfr4', 'ECs', 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 3, 1, 0.903185, 2.709555, 0.03154869388788047, 0.03154869388788046, 0.06796437018833335, 0.611679331695, 0.004631970941917092, 0.004631970941917091),
    @('FAPs', 'Fgf17', 'Fgfr4', 'FAPs', 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 2, 0.6666666666666666, 0.110372, 0.331116, 0.003855347953955327, 0.003855347953955326, 0.008305456209333333, 0.07474910588399999, 0.0005660411729615451, 0.000566041172961545),
    @('FAPs', 'Fgf17', 'Fgfr4', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 1, 0.3333333333333333, 0.03076233333333333, 0.09228699999999999, 0.001074543352259254, 0.001074543352259254, 0.002314855329222222, 0.020833697963, 0.0001577641724625271, 0.0001577641724625271),
    @('FAPs', 'Fgf17', 'Fgfr4', 'MuSCs', 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 3, 1, 27.52907633333334, 82.58722900000001, 0.9616041035407232, 0.9616041035407231, 2.071553817724556, 18.643984359521, 0.1411824616593694, 0.1411824616593693),
    @('FAPs', 'Fgf17', 'Fgfr4', 'Resolving-Mac', 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 2, 0.6666666666666666, 0.05488933333333334, 0.164668, 0.001917311265181737, 0.001917311265181736, 0.004130404036888889, 0.037173636332, 0.0002814991358594322, 0.0002814991358594321),
    @('MuSCs', 'Fgf17', 'Fgfr4', 'ECs', 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 3, 1, 0.903185, 2.709555, 0.03154869388788047, 0.03154869388788046, 0.23085589237, 2.07770303133, 0.01573350539797597, 0.01573350539797597),
    @('MuSCs', 'Fgf17', 'Fgfr4', 'FAPs', 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 2, 0.6666666666666666, 0.110372, 0.331116, 0.003855347953955327, 0.003855347953955326, 0.028211303944, 0.253901735496, 0.001922683013762854, 0.001922683013762854),
    @('MuSCs', 'Fgf17', 'Fgfr4', 'Inflammatory-Mac', 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 1, 0.3333333333333333, 0.03076233333333333, 0.09228699999999999, 0.001074543352259254, 0.001074543352259254, 0.007862913924666667, 0.070766225322, 0.0005358806197560146, 0.0005358806197560146),
    @('MuSCs', 'Fgf17', 'Fgfr4', 'MuSCs', 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 3, 1, 27.52907633333334, 82.58722900000001, 0.9616041035407232, 0.9616041035407231, 7.036486968952667, 63.32838272057401, 0.4795572015609122, 0.4795572015609122),
    @('MuSCs', 'Fgf17', 'Fgfr4', 'Resolving-Mac', 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 2, 0.6666666666666666, 0.05488933333333334, 0.164668, 0.001917311265181737, 0.001917311265181736, 0.01402982337866667, 0.126268410408, 0.0009561735660925528, 0.0009561735660925528)
)

$nr = $rows.Count
$nc = $rows[0].Count
$arr = New-Object 'object[,]' $nr,$nc
for ($i = 0; $i -lt $nr; $i++) {
    for ($j = 0; $j -lt $nc; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $nr - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $nc))
$targetRange.Value = $arr
